$wb = $excel.ActiveWorkbook
$wsImages = $wb.Worksheets.Item("Images")
$wsAttachments = $wb.Worksheets.Item("Attachments")

# Fix typo in the Attachments example id
$wsAttachments.Range("A2").Value = "example_attachment"

# Fix the Url display text so it matches the actual hyperlink target
$wsAttachments.Range("E2").Value = "http://image.url.com"

# Fix the Type value and match the formatting used on the Images sheet
$wsAttachments.Range("F2").Value = "Image"
$wsImages.Range("F2").Copy()
$wsAttachments.Range("F2").PasteSpecial(-4122)

# Leave the Attachments sheet active/selected, as it was the last sheet touched
$wsAttachments.Activate() | Out-Null
$wsAttachments.Range("G2").Select() | Out-Null
